$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.218.36"
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').Value = "'1.687.66"
$ws.Range('E3').Value = '  +0.60%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'216.21"
$ws.Range('E5').Value = '  +0.50%  '
$ws.Range('D6').Value = "'0.522"
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = "'23.07"
$ws.Range('E8').Value = '  +13.31%  '
$ws.Range('E9').Value = '  +3.78%  '
$ws.Range('D11').Value = "'0.0890"
$ws.Range('E11').Value = '  +0.23%  '
$ws.Range('D12').Value = "'1.924.79"
$ws.Range('E12').Value = '  +0.59%  '
$ws.Range('D13').Value = "'1.713.80"
$ws.Range('E13').Value = '  +2.62%  '
$ws.Range('D14').Value = "'4.19"
$ws.Range('E14').Value = '  +2.35%  '
$ws.Range('D15').Value = "'0.552"
$ws.Range('E15').Value = '  +4.36%  '
$ws.Range('D16').Value = "'67.31"
$ws.Range('E16').Value = '  +2.36%  '
$ws.Range('D17').Value = "'27.212.87"
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = "'238.05"
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('D19').Value = "'8.16"
$ws.Range('E19').Value = '  -0.42%  '
$ws.Range('D20').Value = "'0.0₃0746"
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('E22').Value = '  +2.83%  '
$ws.Range('D23').Value = "'9.67"
$ws.Range('E23').Value = '  +5.11%  '
$ws.Range('E24').Value = '  -3.10%  '
$ws.Range('D25').Value = "'148.33"
$ws.Range('E25').Value = '  +1.60%  '
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = "'16.51"
$ws.Range('E27').Value = '  +2.28%  '
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('E29').Value = '  -0.25%  '
$ws.Range('D30').Value = "'0.0503"
$ws.Range('E31').Value = '  +0.39%  '
$ws.Range('B32').Value = 'Maker'
$ws.Range('C32').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D32').Value = "'1.579.25"
$ws.Range('E32').Value = '  +6.37%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'3.40"
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('E34').Value = '  +2.27%  '
$ws.Range('E35').Value = '  +0.24%  '
$ws.Range('D36').Value = "'0.960"
$ws.Range('E37').Value = '  +3.65%  '
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +4.26%  '
$ws.Range('D41').Value = "'69.63"
$ws.Range('E41').Value = '  +3.09%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').Value = "'5.74"
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('B43').Value = 'PaxDollar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D43').Value = "'1.00"
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('D45').Value = "'1.833.37"
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').Value = "'0.787"
$ws.Range('E46').Value = '  +0.77%  '
$ws.Range('D47').Value = "'91.35"
$ws.Range('E47').Value = '  +1.14%  '
$ws.Range('E48').Value = '  +5.81%  '
$ws.Range('E49').Value = '  +3.25%  '
$ws.Range('D51').Value = "'8.23"
$ws.Range('E51').Value = '  +6.09%  '
